$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 45436

$ws.Range("D28").Value = 1053.098
$ws.Range("D29").Value = 1089.132
$ws.Range("D30").Value = 1108.506
$ws.Range("D31").Value = 1141.764
$ws.Range("D32").Value = 1635.074
$ws.Range("D33").Value = 1576.826
$ws.Range("D34").Value = 2203.18
$ws.Range("D35").Value = 2261.376
